# AOS with Parameters and value in data table
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Header cell (becomes the sole shared string "HowMany")
$ws.Range("A1").Value = "HowMany"

# Data table values below the header
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 6

# The 3 data cells get an (unlocked) box border on right/top/bottom - no left edge
$data = $ws.Range("A2:A4")
$data.Locked = $false
$data.Borders.ColorIndex = 1
$data.Borders.Item(7).LineStyle = -4142

# Column A widens to fit the header text
$ws.Range("A1:A4").EntireColumn.AutoFit()

# Move the selection on the Global sheet to B8, then restore Action1 as the
# active tab (Global itself does not become the active sheet)
$ws.Activate()
$ws.Range("B8").Select()
$wb.Worksheets.Item("Action1").Activate()
